$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp text (09:04 -> 09:34)
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 09:34"

# 2. Singapur (row 29) gets new case numbers, stays in the same rank position
$ws.Range("B29").Value = 22460
$ws.Range("C29").Value = 753
$ws.Range("D29").Value = 2040
$ws.Range("E29").Value = 20400
$ws.Range("F29").Value = 22
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 20

# 3. Armenia overtakes Oman (rows 67/68 swap rank, Armenia gets updated data)
$ws.Range("A67").Value = "Armenia"
$ws.Range("B67").Value = 3175
$ws.Range("C67").Value = 146
$ws.Range("D67").Value = 1267
$ws.Range("E67").Value = 1864
$ws.Range("F67").Value = 10
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 44

$ws.Range("A68").Value = "Oman"
$ws.Range("B68").Value = 3112
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 1025
$ws.Range("E68").Value = 2071
$ws.Range("F68").Value = 17
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 16

# 4. Letonia overtakes Somalia (rows 97/98 swap rank, Letonia gets updated data)
$ws.Range("A97").Value = "Letonia"
$ws.Range("B97").Value = 930
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 464
$ws.Range("E97").Value = 448
$ws.Range("F97").Value = 2
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 18

$ws.Range("A98").Value = "Somalia"
$ws.Range("B98").Value = 928
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 106
$ws.Range("E98").Value = 778
$ws.Range("F98").Value = 2
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 44

# 5. Uganda overtakes Nepal (rows 160/161 swap rank, Uganda gets updated data)
$ws.Range("A160").Value = "Uganda"
$ws.Range("B160").Value = 114
$ws.Range("C160").Value = 13
$ws.Range("D160").Value = 55
$ws.Range("E160").Value = 59
$ws.Range("F160").Value = 0
$ws.Range("G160").Value = 0
$ws.Range("H160").Value = 0

$ws.Range("A161").Value = "Nepal"
$ws.Range("B161").Value = 109
$ws.Range("C161").Value = 7
$ws.Range("D161").Value = 30
$ws.Range("E161").Value = 79
$ws.Range("F161").Value = 0
$ws.Range("G161").Value = 0
$ws.Range("H161").Value = 0
